# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Map of row number -> new F-column value (applies identically to both sheets).
$updates = @{
    3  = 1018
    4  = 261
    5  = 1393
    6  = 8483
    7  = 57
    11 = 145
    12 = 3415
    14 = 342
    15 = 57
    16 = 951
    17 = 141
    18 = 1093
    19 = 293
    20 = 160
    21 = 2045
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
